$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.179.11'
$ws.Range("E2").Value = '  +2.48%  '

$ws.Range("D3").Value = '2.294.01'
$ws.Range("E3").Value = '  +3.58%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.97'
$ws.Range("E5").Value = '  +0.43%  '

$ws.Range("E6").Value = '  +2.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.52'
$ws.Range("E7").Value = '  +7.87%  '

$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.644'
$ws.Range("E9").Value = '  +3.72%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.07'
$ws.Range("E10").Value = '  -0.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0984'
$ws.Range("E11").Value = '  +4.99%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.04'
$ws.Range("E12").Value = '  -0.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.42'
$ws.Range("E13").Value = '  +4.78%  '

$ws.Range("E14").Value = '  +1.55%  '

$ws.Range("D15").Value = '2.636.72'
$ws.Range("E15").Value = '  +3.30%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.28'
$ws.Range("E16").Value = '  +5.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.876'
$ws.Range("E17").Value = '  +0.86%  '

$ws.Range("D18").Value = '2.294.84'
$ws.Range("E18").Value = '  +3.81%  '

$ws.Range("D19").Value = '43.057.37'
$ws.Range("E19").Value = '  +2.42%  '

$ws.Range("E20").Value = '  +4.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.31'
$ws.Range("E21").Value = '  +2.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.41'
$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.86'
$ws.Range("E23").Value = '  +1.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.23'
$ws.Range("E24").Value = '  +10.18%  '

$ws.Range("E25").Value = '  +0.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.55'
$ws.Range("E26").Value = '  +3.22%  '

$ws.Range("E27").Value = '  -0.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.42'
$ws.Range("E28").Value = '  +0.76%  '

$ws.Range("E29").Value = '  -1.73%  '

$ws.Range("E30").Value = '  -0.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.05'
$ws.Range("E31").Value = '  +0.42%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.05'
$ws.Range("E32").Value = '  +2.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.44'
$ws.Range("E33").Value = '  +6.51%  '

$ws.Range("E34").Value = '  +5.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0816'
$ws.Range("E35").Value = '  +4.89%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.81'
$ws.Range("E36").Value = '  +20.17%  '

$ws.Range("E37").Value = '  +3.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.59'
$ws.Range("E38").Value = '  +12.10%  '

$ws.Range("E39").Value = '  +3.78%  '

$ws.Range("E40").Value = '  -2.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.52'
$ws.Range("E41").Value = '  +20.84%  '

$ws.Range("E42").Value = '  +5.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.96'
$ws.Range("E43").Value = '  +4.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.220'
$ws.Range("E44").Value = '  +12.47%  '

$ws.Range("B45").Value = 'MultiversX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.99'
$ws.Range("E45").Value = '  +0.73%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.11'
$ws.Range("E46").Value = '  +6.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.94'
$ws.Range("E47").Value = '  -2.13%  '

$ws.Range("E48").Value = '  +2.68%  '

$ws.Range("E49").Value = '  -0.02%  '

$ws.Range("E50").Value = '  +2.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '98.89'
$ws.Range("E51").Value = '  +6.40%  '
